# Expand M05 in checklist. (#86)
#
# Inserts 5 new rows (31-35) holding the 5 Scrum sub-bullets that were
# previously only mentioned in the M05 comment text, pushing every row
# from the old row 31 onward down by 5. Also fixes a few bits of text
# (colon -> "met"/semicolon, trailing period, version string) that were
# touched in the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the intro banner text (version bump).
# ---------------------------------------------------------------------
$a1 = $ws.Cells.Item(1,1)
$a1.Value2 = $a1.Value().Replace("versie 1.1.515, 30-07-2018.", "versie 1.1.710, 31-07-2018.")

# ---------------------------------------------------------------------
# 2. Detach the comments that live on/after row 31 -- they need to move
#    down by 5 rows once the new rows are inserted. Capture their text
#    (and author text box) now, before the geometry changes under us.
# ---------------------------------------------------------------------
$rowsToShift = @(31,32,41,42,43,44,45,47,48,49,50,51,52,62,66,67,68,69,70)
$savedComments = @{}
foreach ($r in $rowsToShift) {
    $cell = $ws.Cells.Item($r,2)
    if ($cell.Comment -ne $null) {
        $savedComments[$r] = $cell.Comment.Text()
        $cell.Comment.Delete()
    }
}

# ---------------------------------------------------------------------
# 3. Insert the 5 new rows right before the old row 31 (M06 and
#    everything below shifts down to make room).
# ---------------------------------------------------------------------
$ws.Rows("31:35").Insert()

# Copy the formatting of a normal "sub-item" row (now at row 38, formerly
# row 33) onto the freshly inserted blank rows so they pick up the same
# style (s=6/s=7) as every other bullet row instead of inheriting the
# M05 header style from row 30.
$ws.Range("A38:D38").Copy()
$ws.Range("A31:D35").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Fill in the text of the 5 new bullet rows (B31:B35), promoted out
#    of the M05 comment text into their own checklist lines.
# ---------------------------------------------------------------------
$ws.Cells.Item(31,2).Value2 = "1. Scrum team bestaand uit product owner, ontwikkelteam en Scrum master"
$ws.Cells.Item(32,2).Value2 = "2. Proces met daily scrum, sprints, sprint planning, sprint review, sprint refinement"
$ws.Cells.Item(33,2).Value2 = "3. Definition of Done"
$ws.Cells.Item(34,2).Value2 = "4. Definition of Ready"
$ws.Cells.Item(35,2).Value2 = "5. Product backlog"

# ---------------------------------------------------------------------
# 5. Re-attach the shifted comments 5 rows further down. Two of them
#    also get a small wording fix while we're at it.
# ---------------------------------------------------------------------
foreach ($r in $savedComments.Keys) {
    $text = $savedComments[$r]
    if ($r -eq 42) {
        $text = $text.Replace("in geval van acceptatie: waarom", "in geval van acceptatie; waarom")
        $text = $text.Replace("in geval van verbeteractie: planning", "in geval van verbeteractie; planning")
    }
    $newRow = $r + 5
    $ws.Cells.Item($newRow, 2).AddComment($text) | Out-Null
}

# Fix the wording of the comment that stays on row 30 (M05 itself).
$c30 = $ws.Cells.Item(30,2)
$txt30 = $c30.Comment.Text().Replace("- Proces: daily scrum", "- Proces met daily scrum")
$c30.Comment.Delete()
$c30.AddComment($txt30) | Out-Null

# ---------------------------------------------------------------------
# 6. Drop the trailing period from the "8. oplevering ..." bullet
#    (now on row 45, formerly row 40).
# ---------------------------------------------------------------------
$c45 = $ws.Cells.Item(45,2)
$c45.Value2 = $c45.Value().TrimEnd('.')

# ---------------------------------------------------------------------
# 7. The conditional formatting range doesn't auto-grow with the
#    inserted rows (unlike data validation / merged cells), so extend
#    it by hand to cover the new bottom of the sheet (row 77).
# ---------------------------------------------------------------------
$fcs = $ws.Range("C3:C72").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("C3:C77"))
}
